# Femacal de La Calera - Pepino dulce: weekly price update.
# A new week's worth of price records is inserted at the top of the data
# (rows 11-12), pushing the existing rows 11-26 down to rows 13-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 11, shifting rows 11:26 -> 13:28.
$ws.Rows("11:12").Insert()

# New record 1 (row 11): Primera quality, new week.
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "Femacal de La Calera"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 45134
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 100112043
$ws.Range("G11").Value = "Pepino dulce"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 68
$ws.Range("K11").Value = 20000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 20000
$ws.Range("N11").Value = "`$/caja 15 kilos"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 1333
$ws.Range("Q11").Value = 15
$ws.Range("R11").Value = "Hortaliza"

# New record 2 (row 12): Segunda quality, same new week.
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "Femacal de La Calera"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 45134
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 100112043
$ws.Range("G12").Value = "Pepino dulce"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 15000
$ws.Range("N12").Value = "`$/caja 15 kilos"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 1000
$ws.Range("Q12").Value = 15
$ws.Range("R12").Value = "Hortaliza"
